$d = $word.ActiveDocument

# --- 1. "Processador :" -> "Processador:" + " " (split into two runs) ---
# Locate the unique occurrence "Processador :" (space before colon) inside the
# "Computador 2:" block. The other occurrence ("Processador: ") already has the
# colon immediately after the word and a trailing space, so it is left untouched.
$rng = $d.Content
$found = $rng.Find.Execute("Processador :", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $matchEnd = $rng.End
    # The matched text is "Processador :" - the space sits right before the colon,
    # i.e. at offset (End-2, End-1).
    $spaceRng = $d.Range($matchEnd - 2, $matchEnd - 1)
    $spaceRng.Text = ""

    # Insert a new space right after the colon (colon is now the last character,
    # ending at $matchEnd - 1).
    $insertPos = $matchEnd - 1
    $insertRng = $d.Range($insertPos, $insertPos)
    $insertRng.Text = " "

    # Toggling a character formatting property forces the engine to keep this
    # newly typed space as its own run instead of silently re-merging it with
    # the neighbouring run that happens to share the same formatting.
    $newSpaceRng = $d.Range($insertPos, $insertPos + 1)
    $newSpaceRng.Bold = 1
    $newSpaceRng.Bold = 0
}

# --- 2. Move the page-number text box further left in the header ---
# margin-left:-49pt -> margin-left:-73.5pt
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
for ($i = 1; $i -le $hdr.Shapes.Count; $i++) {
    $shp = $hdr.Shapes.Item($i)
    if ($shp.Name -eq "Caixa de Texto 22") {
        $shp.Left = -73.5
    }
}
